{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Commit: \"Finish posturize algorithm - no testing\"\n// Adds a new \"Vision (Shaheen)\" section (heading + notes + two reference\n// hyperlinks about RGB<->HSV colour-space conversion for a posturize\n// algorithm) plus a trailing note about an MM slave, all inserted before\n// the document's pre-existing trailing empty paragraph.\n\nconst RAPIDTABLES_URL = \"https://www.rapidtables.com/convert/color/rgb-to-hsv.html\";\nconst HAD2KNOW_URL = \"https://www.had2know.org/technology/hsv-rgb-conversion-formula-calculator.html\";\n\nconst ooxml = `<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Heading2\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>Vision (Shaheen)</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>Convert to HSV, easy colour manipulation</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:hyperlink r:id=\"rIdHyperlink1\" w:history=\"1\">\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t xml:space=\"preserve\">RGB to HSV conversion | </w:t>\n              </w:r>\n              <w:proofErr w:type=\"spellStart\"/>\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t>color</w:t>\n              </w:r>\n              <w:proofErr w:type=\"spellEnd\"/>\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t xml:space=\"preserve\"> conversion (rapidtables.com)</w:t>\n              </w:r>\n            </w:hyperlink>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">Separate H into discrete bands to posturize </w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>Round to nearest band and convert back to RGB</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:hyperlink r:id=\"rIdHyperlink2\" w:history=\"1\">\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t>RGB to HSV, HSV to RGB Conversion Calculator (had2know.org)</w:t>\n              </w:r>\n            </w:hyperlink>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">Use MM slave to allow variable access to NIOS 2 software instead of recompiling to tweak variables </w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/_rels/document.xml.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"256\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rIdHyperlink1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink\" Target=\"${RAPIDTABLES_URL}\" TargetMode=\"External\"/>\n        <Relationship Id=\"rIdHyperlink2\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink\" Target=\"${HAD2KNOW_URL}\" TargetMode=\"External\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n// Insert the whole block at the very start of the body, ahead of the\n// existing (lone) trailing empty paragraph, which must stay last.\nconst startRange = context.document.body.getRange(\"Start\");\nstartRange.insertOoxml(ooxml, \"Start\");\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is the open document.\n#\n# Commit: \"Finish posturize algorithm - no testing\"\n# Adds a new \"Vision (Shaheen)\" section (heading + notes + two reference\n# hyperlinks about RGB<->HSV colour-space conversion for a posturize\n# algorithm) plus a trailing note about an MM slave, all inserted before\n# the document's pre-existing trailing empty paragraph (which is kept,\n# unchanged, as the very last paragraph of the body).\n\n$d = $word.ActiveDocument\n\n$RAPIDTABLES_URL = \"https://www.rapidtables.com/convert/color/rgb-to-hsv.html\"\n$HAD2KNOW_URL = \"https://www.had2know.org/technology/hsv-rgb-conversion-formula-calculator.html\"\n\n# A true zero-length range at the very start of the story (NOT\n# Content.Collapse(), which this host treats as \"whole body\") so the\n# insert lands *before* the existing (lone) trailing empty paragraph\n# and leaves it completely untouched as the last paragraph of the body.\n$target = $d.Range(0, 0)\n\n$ooxml = @\"\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\" pkg:padding=\"512\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:r=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Heading2\"/>\n            </w:pPr>\n            <w:r>\n              <w:t>Vision (Shaheen)</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>Convert to HSV, easy colour manipulation</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:hyperlink r:id=\"rIdHyperlink1\" w:history=\"1\">\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t xml:space=\"preserve\">RGB to HSV conversion | </w:t>\n              </w:r>\n              <w:proofErr w:type=\"spellStart\"/>\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t>color</w:t>\n              </w:r>\n              <w:proofErr w:type=\"spellEnd\"/>\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t xml:space=\"preserve\"> conversion (rapidtables.com)</w:t>\n              </w:r>\n            </w:hyperlink>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">Separate H into discrete bands to posturize </w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:r>\n              <w:t>Round to nearest band and convert back to RGB</w:t>\n            </w:r>\n          </w:p>\n          <w:p>\n            <w:hyperlink r:id=\"rIdHyperlink2\" w:history=\"1\">\n              <w:r>\n                <w:rPr>\n                  <w:color w:val=\"0000FF\"/>\n                  <w:u w:val=\"single\"/>\n                </w:rPr>\n                <w:t>RGB to HSV, HSV to RGB Conversion Calculator (had2know.org)</w:t>\n              </w:r>\n            </w:hyperlink>\n          </w:p>\n          <w:p/>\n          <w:p>\n            <w:r>\n              <w:t xml:space=\"preserve\">Use MM slave to allow variable access to NIOS 2 software instead of recompiling to tweak variables </w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n  <pkg:part pkg:name=\"/word/_rels/document.xml.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\" pkg:padding=\"256\">\n    <pkg:xmlData>\n      <Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">\n        <Relationship Id=\"rIdHyperlink1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink\" Target=\"$RAPIDTABLES_URL\" TargetMode=\"External\"/>\n        <Relationship Id=\"rIdHyperlink2\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/hyperlink\" Target=\"$HAD2KNOW_URL\" TargetMode=\"External\"/>\n      </Relationships>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$target.InsertXML($ooxml)\n"}
